$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1299
$ws.Range("F4").Value = 1031
$ws.Range("F5").Value = 984
$ws.Range("F6").Value = 1758
$ws.Range("F7").Value = 477
$ws.Range("F8").Value = 1173
$ws.Range("F11").Value = 121
$ws.Range("F12").Value = 272
$ws.Range("F14").Value = 82
$ws.Range("F15").Value = 656
$ws.Range("F16").Value = 144
$ws.Range("F17").Value = 94
$ws.Range("F21").Value = 124
$ws.Range("F22").Value = 654
$ws.Range("F23").Value = 21
$ws.Range("F24").Value = 635
$ws.Range("F25").Value = 146
$ws.Range("F27").Value = 857
$ws.Range("F29").Value = 138
$ws.Range("F30").Value = 34
$ws.Range("F31").Value = 259
$ws.Range("F32").Value = 10
$ws.Range("F33").Value = 13
$ws.Range("F34").Value = 404

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 315
$ws.Range("F10").Value = 616
$ws.Range("F12").Value = 21

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1299
$ws.Range("F5").Value = 1031
$ws.Range("F6").Value = 984
$ws.Range("F7").Value = 1759
$ws.Range("F8").Value = 477
$ws.Range("F9").Value = 1173
$ws.Range("F13").Value = 121
$ws.Range("F14").Value = 272
$ws.Range("F16").Value = 82
$ws.Range("F17").Value = 656
$ws.Range("F18").Value = 144
$ws.Range("F19").Value = 94
$ws.Range("F22").Value = 315
$ws.Range("F29").Value = 124
$ws.Range("F30").Value = 654
$ws.Range("F31").Value = 21
$ws.Range("F32").Value = 635
$ws.Range("F33").Value = 146
$ws.Range("F35").Value = 857
$ws.Range("F39").Value = 138
$ws.Range("F40").Value = 34
$ws.Range("F41").Value = 259
$ws.Range("F42").Value = 616
$ws.Range("F45").Value = 10
$ws.Range("F46").Value = 13
$ws.Range("F47").Value = 21
$ws.Range("F48").Value = 404
